# Update crypto price/volume data (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain text before writing new price strings,
# so values like "1.001" / "0.7194" are not auto-coerced to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.216.54"
$ws.Range("E2").Value = "  -0.32%  "
$ws.Range("D3").Value = "1.865.19"
$ws.Range("E3").Value = "  -0.04%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "0.7194"
$ws.Range("E5").Value = "  +2.15%  "
$ws.Range("D6").Value = "240.81"
$ws.Range("E6").Value = "  +1.06%  "
$ws.Range("D7").Value = "1.002"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "0.07763"
$ws.Range("E8").Value = "  +0.41%  "
$ws.Range("D9").Value = "0.3076"
$ws.Range("E9").Value = "  +0.73%  "
$ws.Range("D10").Value = "24.97"
$ws.Range("E10").Value = "  +0.06%  "
$ws.Range("D11").Value = "0.08257"
$ws.Range("E11").Value = "  +0.88%  "
$ws.Range("D12").Value = "1.885.25"
$ws.Range("E12").Value = "  +0.40%  "
$ws.Range("D13").Value = "0.7170"
$ws.Range("E13").Value = "  -0.01%  "
$ws.Range("D14").Value = "5.218"
$ws.Range("E14").Value = "  -0.32%  "
$ws.Range("D15").Value = "90.27"
$ws.Range("E15").Value = "  +1.28%  "
$ws.Range("D16").Value = "29.233.21"
$ws.Range("E16").Value = "  -0.58%  "
$ws.Range("D17").Value = "5.831"
$ws.Range("E17").Value = "  +0.51%  "
$ws.Range("D18").Value = "243.52"
$ws.Range("E18").Value = "  +1.10%  "
$ws.Range("D19").Value = "0.000007792"
$ws.Range("E19").Value = "  -0.15%  "
$ws.Range("D20").Value = "2.121.03"
$ws.Range("E20").Value = "  +0.12%  "
$ws.Range("D21").Value = "13.14"
$ws.Range("E21").Value = "  -0.58%  "
$ws.Range("E22").Value = "  -0.11%  "
$ws.Range("D23").Value = "7.955"
$ws.Range("E23").Value = "  +4.33%  "
$ws.Range("D24").Value = "1.002"
$ws.Range("E24").Value = "  -0.14%  "
$ws.Range("D25").Value = "0.1591"
$ws.Range("E25").Value = "  +9.96%  "
$ws.Range("D26").Value = "162.34"
$ws.Range("E26").Value = "  +0.18%  "
$ws.Range("D27").Value = "8.913"
$ws.Range("E27").Value = "  -0.24%  "
$ws.Range("D28").Value = "18.18"
$ws.Range("E28").Value = "  +0.24%  "
$ws.Range("D29").Value = "1.495"
$ws.Range("E29").Value = "  -0.52%  "
$ws.Range("D30").Value = "1.309"
$ws.Range("E30").Value = "  -4.40%  "
$ws.Range("D31").Value = "4.354"
$ws.Range("E31").Value = "  +1.14%  "
$ws.Range("D32").Value = "4.086"
$ws.Range("E32").Value = "  +1.12%  "
$ws.Range("D33").Value = "0.05184"
$ws.Range("E33").Value = "  -0.48%  "
$ws.Range("D34").Value = "1.916"
$ws.Range("E34").Value = "  -0.14%  "
$ws.Range("D35").Value = "1.176"
$ws.Range("E35").Value = "  -0.92%  "
$ws.Range("D36").Value = "0.7283"
$ws.Range("E36").Value = "  +1.83%  "
$ws.Range("D37").Value = "2.678"
$ws.Range("E37").Value = "  -0.11%  "
$ws.Range("D38").Value = "0.01849"
$ws.Range("E38").Value = "  +0.03%  "
$ws.Range("D39").Value = "2.696"
$ws.Range("E39").Value = "  -0.11%  "
$ws.Range("D40").Value = "1.158.34"
$ws.Range("E40").Value = "  -1.42%  "
$ws.Range("D41").Value = "0.9020"
$ws.Range("E41").Value = "  -1.47%  "
$ws.Range("D42").Value = "6.109"
$ws.Range("E42").Value = "  +2.06%  "
$ws.Range("D43").Value = "72.27"
$ws.Range("E43").Value = "  +1.67%  "
$ws.Range("D44").Value = "1.002"
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("D45").Value = "101.74"
$ws.Range("E45").Value = "  -0.27%  "
$ws.Range("D46").Value = "2.018.39"
$ws.Range("E46").Value = "  +0.14%  "
$ws.Range("D47").Value = "0.5273"
$ws.Range("E47").Value = "  -1.73%  "
$ws.Range("D48").Value = "1.763"
$ws.Range("E48").Value = "  +0.80%  "
$ws.Range("D49").Value = "9.287"
$ws.Range("E49").Value = "  +0.96%  "
$ws.Range("D50").Value = "2.869"
$ws.Range("E50").Value = "  +2.07%  "
$ws.Range("D51").Value = "0.9999"
$ws.Range("E51").Value = "  -0.44%  "

# Restore default (unstyled) formatting on column D now that the text
# values are committed, matching the original workbook styling.
$ws.Range("D2:D51").ClearFormats()
